# Update the "K" column (column G) values for rows 2-19 in the active sheet.
# This reflects the regen of save_data to use K instead of Strike# (commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 3
    11 = 3
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    19 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
